$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty D33 cell (was present in the source, removed in target)
$ws.Cells.Item(33, 4).ClearContents()

# Row 34
$ws.Cells.Item(34, 1).Value = "2025-09-16 08:47:58"
$ws.Cells.Item(34, 2).Value = "demo"
$ws.Cells.Item(34, 3).Value = "Ford"
$ws.Cells.Item(34, 4).Value = "Truecue iEngineering"
$ws.Cells.Item(34, 5).Value = "sophie.jones@ienguast.com"
$ws.Cells.Item(34, 6).Value = "AT"
$ws.Cells.Item(34, 7).Value = "Austria"

# Row 35
$ws.Cells.Item(35, 1).Value = "2025-09-16 08:48:29"
$ws.Cells.Item(35, 2).Value = "demo"
$ws.Cells.Item(35, 3).Value = "Bhutekar"
$ws.Cells.Item(35, 5).Value = "kushankur@iengaust.com.au"
$ws.Cells.Item(35, 6).Value = "BD"
$ws.Cells.Item(35, 7).Value = "Bangladesh"

# Row 36
$ws.Cells.Item(36, 1).Value = "2025-09-16 08:52:34"
$ws.Cells.Item(36, 2).Value = "demo"
$ws.Cells.Item(36, 3).Value = "Ford"
$ws.Cells.Item(36, 5).Value = "kushankur@iengaust.com.au"
$ws.Cells.Item(36, 6).Value = "AD"
$ws.Cells.Item(36, 7).Value = "Andorra"

# Row 37
$ws.Cells.Item(37, 1).Value = "2025-09-16 09:39:46"
$ws.Cells.Item(37, 2).Value = "demo"
$ws.Cells.Item(37, 3).Value = "Ford"
$ws.Cells.Item(37, 4).Value = "Truecue iEngineering"
$ws.Cells.Item(37, 5).Value = "kushankur@iengaust.com.au"
$ws.Cells.Item(37, 6).Value = "IN"
$ws.Cells.Item(37, 7).Value = "India"
$ws.Cells.Item(37, 8).Value = "'+919168627258"

# Row 38
$ws.Cells.Item(38, 1).Value = "2025-09-16 10:00:54"
$ws.Cells.Item(38, 2).Value = "demo"
$ws.Cells.Item(38, 3).Value = "Ford"
$ws.Cells.Item(38, 5).Value = "mayur@iengaust.com.au"
$ws.Cells.Item(38, 6).Value = "IN"
$ws.Cells.Item(38, 7).Value = "India"
$ws.Cells.Item(38, 8).Value = "'+919168627258"

# Row 39
$ws.Cells.Item(39, 1).Value = "2025-09-19 15:27:02"
$ws.Cells.Item(39, 2).Value = "Diksha"
$ws.Cells.Item(39, 3).Value = "diksha@iengaust.com.au"
$ws.Cells.Item(39, 4).Value = "'+918956972428"
$ws.Cells.Item(39, 5).Value = "IN"
$ws.Cells.Item(39, 6).Value = "India"
$ws.Cells.Item(39, 7).Value = "'+91"
$ws.Cells.Item(39, 8).Value = "IoT / Sensors"
$ws.Cells.Item(39, 9).Value = "ddddd"
$ws.Cells.Item(39, 10).Value = "127.0.0.1"
$ws.Cells.Item(39, 11).Value = "http://127.0.0.1:8000/neplan-electricity/"

# Row 40
$ws.Cells.Item(40, 1).Value = "2025-09-19 15:47:29"
$ws.Cells.Item(40, 2).Value = "Diksha"
$ws.Cells.Item(40, 3).Value = "diksha@iengaust.com.au"
$ws.Cells.Item(40, 4).Value = "'+919168627258"
$ws.Cells.Item(40, 5).Value = "IN"
$ws.Cells.Item(40, 6).Value = "India"
$ws.Cells.Item(40, 7).Value = "'+91"
$ws.Cells.Item(40, 8).Value = "NEPLAN Gas, Water and Heating"
$ws.Cells.Item(40, 9).Value = "hi"
$ws.Cells.Item(40, 10).Value = "127.0.0.1"
$ws.Cells.Item(40, 11).Value = "http://127.0.0.1:8000/neplan-electricity/"

# Row 41
$ws.Cells.Item(41, 1).Value = "2025-09-19 15:47:57"
$ws.Cells.Item(41, 2).Value = "Diksha"
$ws.Cells.Item(41, 3).Value = "dnaiker@iengaust.com.au"
$ws.Cells.Item(41, 4).Value = "'+919168627258"
$ws.Cells.Item(41, 5).Value = "IN"
$ws.Cells.Item(41, 6).Value = "India"
$ws.Cells.Item(41, 7).Value = "'+91"
$ws.Cells.Item(41, 10).Value = "127.0.0.1"
$ws.Cells.Item(41, 11).Value = "http://127.0.0.1:8000/neplan-electricity/"

# Row 42
$ws.Cells.Item(42, 1).Value = "2025-09-19 16:08:36"
$ws.Cells.Item(42, 2).Value = "iksha"
$ws.Cells.Item(42, 3).Value = "diksha@iengaust.com.au"
$ws.Cells.Item(42, 4).Value = "'+919168627258"
$ws.Cells.Item(42, 5).Value = "IN"
$ws.Cells.Item(42, 6).Value = "India"
$ws.Cells.Item(42, 7).Value = "'+91"
$ws.Cells.Item(42, 8).Value = "NEPLAN Electricity"
$ws.Cells.Item(42, 10).Value = "127.0.0.1"
$ws.Cells.Item(42, 11).Value = "http://127.0.0.1:8000/neplan-electricity/"

# Row 43
$ws.Cells.Item(43, 1).Value = "2025-09-19 16:15:02"
$ws.Cells.Item(43, 2).Value = "iksha"
$ws.Cells.Item(43, 3).Value = "diksha@iengaust.com.au"
$ws.Cells.Item(43, 4).Value = "'+918956972428"
$ws.Cells.Item(43, 5).Value = "IN"
$ws.Cells.Item(43, 6).Value = "India"
$ws.Cells.Item(43, 7).Value = "'+91"
$ws.Cells.Item(43, 8).Value = "NEPLAN Electricity"
$ws.Cells.Item(43, 10).Value = "127.0.0.1"
$ws.Cells.Item(43, 11).Value = "http://127.0.0.1:8000/neplan-electricity/"

# Row 44
$ws.Cells.Item(44, 1).Value = "2025-09-19 16:16:40"
$ws.Cells.Item(44, 2).Value = "iksha"
$ws.Cells.Item(44, 3).Value = "diksha@iengaust.com.au"
$ws.Cells.Item(44, 4).Value = "'+919168627258"
$ws.Cells.Item(44, 5).Value = "IN"
$ws.Cells.Item(44, 6).Value = "India"
$ws.Cells.Item(44, 7).Value = "'+91"
$ws.Cells.Item(44, 8).Value = "NEPLAN Electricity"
$ws.Cells.Item(44, 9).Value = "gg"
$ws.Cells.Item(44, 10).Value = "127.0.0.1"
$ws.Cells.Item(44, 11).Value = "http://127.0.0.1:8000/neplan-electricity/"

# Row 45
$ws.Cells.Item(45, 1).Value = "2025-09-24 08:26:50"
$ws.Cells.Item(45, 2).Value = "Mayur Mane"
$ws.Cells.Item(45, 3).Value = "mayur@iengaust.com.au"
$ws.Cells.Item(45, 4).Value = "'+9189897889785"
$ws.Cells.Item(45, 5).Value = "IN"
$ws.Cells.Item(45, 6).Value = "India"
$ws.Cells.Item(45, 7).Value = "'+91"
$ws.Cells.Item(45, 8).Value = "NEPLAN Electricity"
$ws.Cells.Item(45, 10).Value = "192.168.1.53"
$ws.Cells.Item(45, 11).Value = "http://192.168.1.58:8000/contact/"

# Row 46
$ws.Cells.Item(46, 1).Value = "2025-09-24 08:27:56"
$ws.Cells.Item(46, 2).Value = "Mayur Mane"
$ws.Cells.Item(46, 3).Value = "mayur@iengaust.com.au"
$ws.Cells.Item(46, 4).Value = "'+9189897889785"
$ws.Cells.Item(46, 5).Value = "IN"
$ws.Cells.Item(46, 6).Value = "India"
$ws.Cells.Item(46, 7).Value = "'+91"
$ws.Cells.Item(46, 8).Value = "NEPLAN Gas, Water and Heating"
$ws.Cells.Item(46, 10).Value = "192.168.1.53"
$ws.Cells.Item(46, 11).Value = "http://192.168.1.58:8000/contact/"

# Row 47
$ws.Cells.Item(47, 1).Value = "2025-09-24 08:30:11"
$ws.Cells.Item(47, 2).Value = "Mayur Mane"
$ws.Cells.Item(47, 3).Value = "mayur@iengaust.com.au"
$ws.Cells.Item(47, 4).Value = "'+9189897889785"
$ws.Cells.Item(47, 5).Value = "IN"
$ws.Cells.Item(47, 6).Value = "India"
$ws.Cells.Item(47, 7).Value = "'+91"
$ws.Cells.Item(47, 8).Value = "NEPLAN Electricity"
$ws.Cells.Item(47, 10).Value = "192.168.1.53"
$ws.Cells.Item(47, 11).Value = "http://192.168.1.58:8000/contact/"

# Row 48
$ws.Cells.Item(48, 1).Value = "2025-09-24 08:47:17"
$ws.Cells.Item(48, 2).Value = "Mayur Mane"
$ws.Cells.Item(48, 3).Value = "mayur@iengaust.com.au"
$ws.Cells.Item(48, 4).Value = "'+9189897889785"
$ws.Cells.Item(48, 5).Value = "IN"
$ws.Cells.Item(48, 6).Value = "India"
$ws.Cells.Item(48, 7).Value = "'+91"
$ws.Cells.Item(48, 8).Value = "NEPLAN Gas, Water and Heating"
$ws.Cells.Item(48, 10).Value = "192.168.1.53"
$ws.Cells.Item(48, 11).Value = "http://192.168.1.58:8000/contact/"

# Row 49
$ws.Cells.Item(49, 1).Value = "2025-09-24 09:46:01"
$ws.Cells.Item(49, 2).Value = "Mayur Mane"
$ws.Cells.Item(49, 3).Value = "mayur@iengaust.com.au"
$ws.Cells.Item(49, 4).Value = "'+9189897889785"
$ws.Cells.Item(49, 5).Value = "IN"
$ws.Cells.Item(49, 6).Value = "India"
$ws.Cells.Item(49, 7).Value = "'+91"
$ws.Cells.Item(49, 8).Value = "NEPLAN Electricity"
$ws.Cells.Item(49, 10).Value = "192.168.1.53"
$ws.Cells.Item(49, 11).Value = "http://192.168.1.58:8000/neplan-electricity/"

# Row 50
$ws.Cells.Item(50, 1).Value = "2025-09-30 16:50:24"
$ws.Cells.Item(50, 2).Value = "Diksha"
$ws.Cells.Item(50, 3).Value = "dikshabhutekar@gmail.com"
$ws.Cells.Item(50, 4).Value = "'+918956972428"
$ws.Cells.Item(50, 5).Value = "IN"
$ws.Cells.Item(50, 6).Value = "India"
$ws.Cells.Item(50, 7).Value = "'+91"
$ws.Cells.Item(50, 8).Value = "NEPLAN OEM / SaS Webservice"
$ws.Cells.Item(50, 9).Value = "'"
$ws.Cells.Item(50, 10).Value = "127.0.0.1"
$ws.Cells.Item(50, 11).Value = "http://127.0.0.1:8000/contact/"
